# level_act7d5_st02 - update quotation marks (" -> ') in EN text and
# 의료팀 -> 의료부 in KO text, per commit "update on 20210731 画中人"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# en_US column (C)
$ws.Range("C19").Value = '[name="Diary"]     Ever since he found out that Warfarin is ''Mr. Blood,'' he has been visiting the medical section regularly.
'
$ws.Range("C25").Value = '[name="Ptilopsis"]     This diary must be from the time of his second ''duel'' with Warfarin.
'
$ws.Range("C26").Value = '[name="Aak"]     Hey, ''Mr. Blood,'' how about a little contest of skill?
'
$ws.Range("C37").Value = '[name="Aak"]     Sure, sure. I''m just asking a question. Give me an answer and I''ll leave you to your noble ''work.''
'
$ws.Range("C63").Value = '[name="Ptilopsis"]     Oh, I wrote, ''I did not dislike it'' in this one.
'
$ws.Range("C198").Value = '[name="Ptilopsis"]     Is ''normal'' the correct word here? Is it my life that is abnormal?
'

# ko_KR column (D)
$ws.Range("D17").Value = '[name="일기"]     둘 중에선 아 씨가 인상에 남습니다. 굉장히 젊어 보이지만, 듣기론 용문에서 유명한 무면허 의사였다는 것 같습니다. 의료부에는 들어오지 않을 것 같습니다.
'
$ws.Range("D19").Value = '[name="일기"]     ……그러고 보니, 와파린 씨가 ''미스터 블러드''라는 것을 알게 된 아 씨는 종종 의료부에 찾아오고 있습니다.
'
$ws.Range("D20").Value = '[name="일기"]     하지만 의료부 대다수는 그를 싫어합니다. 그가 우리를 향한 혐오감을 감추지 않았기 때문입니다.
'
$ws.Range("D24").Value = '[name="일기"]     의료부 사람들도 그의 성격에 익숙해져 가고 있습니다…… 적어도 표면상으로는 그래 보입니다.
'
$ws.Range("D62").Value = '[name="일기"]     아 씨가 온 이후로 의료부는 시끌벅적해졌습니다. 물론 유쾌하지 않은 일도 있었지만, 왠지 모르게 싫지는 않았습니다……
'
$ws.Range("D150").Value = '[name="일기"]     저는, 의료부에게 줄 선물을 사는 게 어떨까 생각했습니다.
'
$ws.Range("D159").Value = '[name="일기"]     그때, 저는 두 번째 문제에 직면했습니다. 의료부 동료들이 무엇을 좋아하는지 몰랐기 때문입니다.
'
